$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the old backlog rows (5-17); only the header + 3 task rows remain.
$ws.Range("A5:A17").EntireRow.Delete()

# Rewrite the remaining task list.
$ws.Range("A2").Value = "Autoscout24.de"
$ws.Range("B2").Value = "Artas"
$ws.Range("A3").Value = "Logas parserio"
$ws.Range("B3").Value = "Ignas"
$ws.Range("A4").Value = "Edit disable settings"
$ws.Range("B4").Value = "Ignas"

# Standard/Normal formatting for the data rows (drop the "Good" green style
# and wrap-text formatting that used to mark finished tasks).
$ws.Range("A2:B4").Style = "Normal"

# The "Good" cell style is no longer used anywhere - remove it.
$wb.Styles("Good").Delete()

# Match the selection left behind by the edit.
$null = $ws.Range("A5").Select()
